# Weekly update: a new reporting week's records are prepended to the
# "Pepino ensalada" log, pushing the existing history down by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows above the first data block for this price series
# (row 348) so every existing record shifts down by two rows.
$ws.Rows("348:349").Insert()

# New row 348 - "Primera" quality entry for the new week (2023-03-24).
$ws.Cells.Item(348, 1).Value = 1
$ws.Cells.Item(348, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(348, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(348, 4).Value = 45009
$ws.Cells.Item(348, 5).Value = 15
$ws.Cells.Item(348, 6).Value = 100112043
$ws.Cells.Item(348, 7).Value = "Pepino ensalada"
$ws.Cells.Item(348, 8).Value = "Sin especificar"
$ws.Cells.Item(348, 9).Value = "Primera"
$ws.Cells.Item(348, 10).Value = 350
$ws.Cells.Item(348, 11).Value = 4500
$ws.Cells.Item(348, 12).Value = 5000
$ws.Cells.Item(348, 13).Value = 4786
$ws.Cells.Item(348, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(348, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(348, 16).Value = 68
$ws.Cells.Item(348, 17).Value = 70
$ws.Cells.Item(348, 18).Value = "Hortaliza"

# New row 349 - "Segunda" quality entry for the same new week.
$ws.Cells.Item(349, 1).Value = 1
$ws.Cells.Item(349, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(349, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(349, 4).Value = 45009
$ws.Cells.Item(349, 5).Value = 15
$ws.Cells.Item(349, 6).Value = 100112043
$ws.Cells.Item(349, 7).Value = "Pepino ensalada"
$ws.Cells.Item(349, 8).Value = "Sin especificar"
$ws.Cells.Item(349, 9).Value = "Segunda"
$ws.Cells.Item(349, 10).Value = 180
$ws.Cells.Item(349, 11).Value = 3000
$ws.Cells.Item(349, 12).Value = 3500
$ws.Cells.Item(349, 13).Value = 3222
$ws.Cells.Item(349, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(349, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(349, 16).Value = 32
$ws.Cells.Item(349, 17).Value = 100
$ws.Cells.Item(349, 18).Value = "Hortaliza"
